$d = $word.ActiveDocument

# "Editing Class for admin": append a new sentence about admin starting health
# to the end of the "Class Selection" paragraph's description.
$old = "Knights have high health but deal less damage."
$new = "Knights have high health but deal less damage. Admins are granted 1000 health points at the beginning."

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
